# Insert a new weekly data row above row 174 (Berenjena / Vega Modelo de
# Temuco price sheet). Excel's native row-insert shifts rows 174..207 down
# to 175..208 (and extends the sheet's used range to row 208), then we
# populate the newly-inserted row 174 with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(174).Insert()

$ws.Range("A174").Value = 10
$ws.Range("B174").Value = "Vega Modelo de Temuco"
$ws.Range("C174").Value = "La Araucanía"
$ws.Range("D174").Value = 44522
$ws.Range("E174").Value = 9
$ws.Range("F174").Value = 100112001
$ws.Range("G174").Value = "Berenjena"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 50
$ws.Range("K174").Value = 10000
$ws.Range("L174").Value = 10000
$ws.Range("M174").Value = 10000
$ws.Range("N174").Value = "$/caja 60 unidades"
$ws.Range("O174").Value = "Región de Arica y Parinacota"
$ws.Range("P174").Value = 167
$ws.Range("Q174").Value = 60
$ws.Range("R174").Value = "Hortaliza"
